$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-07-29"

# Update the label for the July row.
$ws.Range("A8").Value = "July (through 07-29)"

# Update the July row figures (columns C:I; B is unchanged).
$ws.Range("C8").Value = 51
$ws.Range("D8").Value = 70
$ws.Range("E8").Value = 68
$ws.Range("F8").Value = 49
$ws.Range("G8").Value = 138
$ws.Range("H8").Value = 141
$ws.Range("I8").Value = 162

# Update the Total row figures (columns C:I; B is unchanged).
$ws.Range("C9").Value = 299
$ws.Range("D9").Value = 460
$ws.Range("E9").Value = 421
$ws.Range("F9").Value = 300
$ws.Range("G9").Value = 610
$ws.Range("H9").Value = 901
$ws.Range("I9").Value = 968
